$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 19:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 829013
$ws.Range("C4").Value = 10269
$ws.Range("E4").Value = 699446
$ws.Range("G4").Value = 829
$ws.Range("H4").Value = 46147

# Row 8 - Alemania
$ws.Range("B8").Value = 149401
$ws.Range("C8").Value = 948
$ws.Range("E8").Value = 44836
$ws.Range("G8").Value = 79
$ws.Range("H8").Value = 5165

# Row 10 - Turquia
$ws.Range("B10").Value = 98674
$ws.Range("C10").Value = 3083
$ws.Range("D10").Value = 16477
$ws.Range("E10").Value = 79821
$ws.Range("F10").Value = 1814
$ws.Range("G10").Value = 117
$ws.Range("H10").Value = 2376

# Row 80 - Cuba
$ws.Range("F80").Value = 16

# Row 93 - Principado de Andorra
$ws.Range("B93").Value = 723
$ws.Range("C93").Value = 6
$ws.Range("D93").Value = 309
$ws.Range("E93").Value = 377
